# "Segundo sprint INSO.xlsx" - Add files via upload
#
# The task description for the sub-task in row 6 (D6) gains " y empleado" at
# the end, and the task description for the sub-task in row 7 (D7) is
# replaced with a new task about the JavaScript/CSS implementation (it used
# to duplicate the PHP/HTML text that is now only in D6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historias de Usuario")

$ws.Range("D6").Value = "Implementación del codigo de la aplicación en PHP, HTML , de la parte del cliente y empleado"
$ws.Range("D7").Value = "Implementación del codigo de la aplicación en JavaScript y el CSS de la parte del cliente y empleado"

# The rows keep their existing (wrapped-text) height, but touching them marks
# the height as an explicit/custom one rather than auto-fit.
$ws.Rows.Item(6).RowHeight = 75
$ws.Rows.Item(7).RowHeight = 75

# Move the active-cell selection on the frozen "bottomRight" pane to D7,
# matching where the author was last working.
$ws.Range("D7").Select()
